# Generate Report for handback
# Updates the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) values on row 5 of the
# "zh-cn" and "de-de" sheets to reflect a newly generated report.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-01-25 14:08:31"
$wsZh.Range("G5").Value = "2016-01-25 14:09:15"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-01-25 14:08:40"
$wsDe.Range("G5").Value = "2016-01-25 14:09:33"
